# Updates cryptos list values per diff: price and 1h-volume-change columns,
# plus reshuffled coin rows 43-51 (names/links/prices/changes).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.988.93'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').Value = '1.641.43'
$ws.Range('E4').Value = '  +0.40%  '
$ws.Range('D5').Value = '''215.81'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('E7').Value = '  +0.41%  '
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('E12').Value = '  +0.52%  '
$ws.Range('D13').Value = '1.645.11'
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').Value = '0.0₃0763'
$ws.Range('E15').Value = '  +0.87%  '
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('D17').Value = '26.021.85'
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('D19').Value = '''193.98'
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('E22').Value = '  -1.05%  '
$ws.Range('E23').Value = '  +4.73%  '
$ws.Range('D24').Value = '''1.79'
$ws.Range('E24').Value = '  -1.47%  '
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('D26').Value = '''142.91'
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('E28').Value = '  +0.60%  '
$ws.Range('E29').Value = '  +0.51%  '
$ws.Range('E30').Value = '  -0.79%  '
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('E32').Value = '  +1.35%  '
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('E34').Value = '  +1.37%  '
$ws.Range('E35').Value = '  +0.27%  '
$ws.Range('D36').Value = '1.129.77'
$ws.Range('E36').Value = '  -0.96%  '
$ws.Range('D37').Value = '''0.540'
$ws.Range('E37').Value = '  -1.12%  '
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('D41').Value = '''98.99'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('B43').Value = 'BabyDogeCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D43').Value = '0.0₆0117'
$ws.Range('E43').Value = '  +4.71%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '''56.56'
$ws.Range('E44').Value = '  +0.41%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '''1.49'
$ws.Range('E45').Value = '  +3.26%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = '''0.0522'
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''7.73'
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '''0.414'
$ws.Range('E48').Value = '  -0.22%  '
$ws.Range('B49').Value = 'USDD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D49').Value = '''1.00'
$ws.Range('E49').Value = '  +0.42%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '''0.0952'
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '''1.19'
$ws.Range('E51').Value = '  +3.78%  '
